# Add "2022-Q1" fund-holdings sheet and refresh the "总计" (total) summary
# sheet with a new leading row for 2022-Q1.
#
# Strategy: the existing "总计" worksheet (sheetId 6) keeps its sheetId/rId
# but is renamed to "2022-Q1" and repopulated with the quarterly fund
# holdings table (same layout as the other quarter sheets). A brand-new
# worksheet named "总计" is appended after it, carrying the old total-table
# data plus one new row for 2022-Q1 at the top.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: repurpose the current "总计" sheet -> "2022-Q1"
# ---------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Item("总计")
$q1Sheet.Name = "2022-Q1"

# Clear any pre-existing content from the old total table.
$q1Sheet.Cells.Clear()

# Borrow cell formatting (fonts/borders/alignment) from the "2021-Q4" sheet,
# which already has the identical 8-column fund-holdings layout, so the new
# sheet reuses the same style records instead of minting new ones.
$styleSource = $wb.Worksheets.Item("2021-Q4")
$styleSource.Range("B1:H1").Copy()
$q1Sheet.Range("B1:H1").PasteSpecial(-4122)
$styleSource.Range("A2:H13").Copy()
$q1Sheet.Range("A2:H13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Data rows (A=index, B=fund code, C=fund name, D..G stored as text,
# H=rank as a number) matching the quarterly-sheet convention.
$q1Data = @(
    @{A=0;  B="519002"; C="华安安信消费混合";            D="80.30"; E="89.00"; F="2.78"; G="2.2323"; H=4},
    @{A=1;  B="001532"; C="华安文体健康主题灵活配置混合"; D="87.24"; E="92.86"; F="2.31"; G="2.0152"; H=7},
    @{A=2;  B="011251"; C="华安聚嘉精选混合A";           D="31.29"; E="89.00"; F="2.92"; G="0.9137"; H=4},
    @{A=3;  B="011128"; C="华安精致生活混合A";           D="33.30"; E="85.22"; F="2.62"; G="0.8725"; H=3},
    @{A=4;  B="011252"; C="华安聚嘉精选混合C";           D="14.89"; E="89.00"; F="2.92"; G="0.4348"; H=4},
    @{A=5;  B="000780"; C="鹏华医疗保健股票";             D="7.67";  E="82.80"; F="3.67"; G="0.2815"; H=5},
    @{A=6;  B="011129"; C="华安精致生活混合C";           D="7.57";  E="85.22"; F="2.62"; G="0.1983"; H=3},
    @{A=7;  B="011471"; C="鹏华致远成长混合A";           D="2.19";  E="61.03"; F="3.17"; G="0.0694"; H=3},
    @{A=8;  B="009956"; C="广发恒誉混合A";               D="4.94";  E="21.40"; F="1.27"; G="0.0627"; H=1},
    @{A=9;  B="005295"; C="诺德天富灵活配置混合";         D="1.21";  E="93.81"; F="1.94"; G="0.0235"; H=2},
    @{A=10; B="011472"; C="鹏华致远成长混合C";           D="0.08";  E="61.03"; F="3.17"; G="0.0025"; H=3},
    @{A=11; B="009957"; C="广发恒誉混合C";               D="0.10";  E="21.40"; F="1.27"; G="0.0013"; H=1}
)

$r = 2
foreach ($row in $q1Data) {
    $q1Sheet.Range("A$r").Value = $row.A
    $q1Sheet.Range("B$r").Value = "'" + $row.B
    $q1Sheet.Range("C$r").Value = $row.C
    $q1Sheet.Range("D$r").Value = "'" + $row.D
    $q1Sheet.Range("E$r").Value = "'" + $row.E
    $q1Sheet.Range("F$r").Value = "'" + $row.F
    $q1Sheet.Range("G$r").Value = "'" + $row.G
    $q1Sheet.Range("H$r").Value = $row.H
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: append a brand-new "总计" sheet after "2022-Q1"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1Sheet)
$totalSheet.Name = "总计"

# Reuse the same header/data styling from the sheet we just repopulated
# (it still carries the s="2" bold-centered style on column B/C/D header
# cells and column A data cells).
$q1Sheet.Range("B1:D1").Copy()
$totalSheet.Range("B1:D1").PasteSpecial(-4122)
$q1Sheet.Range("A2:D7").Copy()
$totalSheet.Range("A2:D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the page-margin settings used throughout this workbook.
$totalSheet.PageSetup.LeftMargin = 0.75 * 72
$totalSheet.PageSetup.RightMargin = 0.75 * 72
$totalSheet.PageSetup.TopMargin = 1 * 72
$totalSheet.PageSetup.BottomMargin = 1 * 72
$totalSheet.PageSetup.HeaderMargin = 0.5 * 72
$totalSheet.PageSetup.FooterMargin = 0.5 * 72

$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$totalData = @(
    @{A=0; B="2022-Q1"; C=12; D=7.11},
    @{A=1; B="2021-Q4"; C=15; D=9.08},
    @{A=2; B="2021-Q3"; C=6;  D=4.14},
    @{A=3; B="2021-Q2"; C=6;  D=0.79},
    @{A=4; B="2021-Q1"; C=3;  D=0.11},
    @{A=5; B="2020-Q4"; C=5;  D=1.68}
)

$r = 2
foreach ($row in $totalData) {
    $totalSheet.Range("A$r").Value = $row.A
    $totalSheet.Range("B$r").Value = $row.B
    $totalSheet.Range("C$r").Value = $row.C
    $totalSheet.Range("D$r").Value = $row.D
    $r = $r + 1
}

# Restore the originally active sheet/tab ("2020-Q4", the first sheet) so
# the workbook-level selection state is unchanged by these edits.
$wb.Worksheets.Item(1).Select()
